# SoIB_summaries.xlsx — re-ran resolve and classify+summarise steps after
# changes to the mapping file. Updates "Range Status", "Species
# qualification" and "High Priority break-up" sheets to reflect the new
# summary numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Range Status" sheet: the range analysis no longer has species-level
# breakdown counts, so column C (Species (perc.)) is cleared for every
# data row, and column B (Species (no.)) collapses to 0 across the board.
# ---------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("B3").Value = 0
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# ---------------------------------------------------------------------
# "Species qualification" sheet: Range Analysis species count drops to 0.
# ---------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# ---------------------------------------------------------------------
# "High Priority break-up" sheet: the "Trend Different" row disappears
# entirely (mapping-file change means no species fall in that bucket any
# more), so IUCN shifts up to row 3, and the remaining rows' counts are
# recomputed against the new totals.
# ---------------------------------------------------------------------
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Rows.Item(3).Delete()

$wsBreak.Range("B2").Value = 3
$wsBreak.Range("C2").Value = 23.1
$wsBreak.Range("D2").Value = 3
$wsBreak.Range("E2").Value = 23.1

$wsBreak.Range("D3").Value = 10
$wsBreak.Range("E3").Value = 76.9
